# Daily attendance processing - 2025-10-31 11:42:35
# Swap the order of the two "Recorded By" entries (column G) for the rows
# identified in the upstream diff. Each change reverses a two-item,
# comma-separated list of recorder names/emails (e.g. "A, B" -> "B, A").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G3" = "System, dnasr281@gmail.com"
    "G4" = "backup@backdoor.com, System"
    "G6" = "System, dnasr281@gmail.com"
    "G7" = "System, admin@admin.com"
    "G10" = "System, dnasr281@gmail.com"
    "G12" = "System, dnasr281@gmail.com"
    "G13" = "System, dnasr281@gmail.com"
    "G14" = "System, dnasr281@gmail.com"
    "G15" = "System, dnasr281@gmail.com"
    "G18" = "System, dnasr281@gmail.com"
    "G19" = "System, dnasr281@gmail.com"
    "G20" = "System, dnasr281@gmail.com"
    "G21" = "System, dnasr281@gmail.com"
    "G22" = "System, dnasr281@gmail.com"
    "G24" = "System, dnasr281@gmail.com"
    "G30" = "System, dnasr281@gmail.com"
    "G31" = "backup@backdoor.com, System"
    "G33" = "System, dnasr281@gmail.com"
    "G34" = "System, admin@admin.com"
    "G37" = "System, dnasr281@gmail.com"
    "G39" = "System, dnasr281@gmail.com"
    "G40" = "System, dnasr281@gmail.com"
    "G41" = "System, dnasr281@gmail.com"
    "G42" = "System, dnasr281@gmail.com"
    "G45" = "System, dnasr281@gmail.com"
    "G46" = "System, dnasr281@gmail.com"
    "G47" = "System, dnasr281@gmail.com"
    "G48" = "System, dnasr281@gmail.com"
    "G49" = "System, dnasr281@gmail.com"
    "G51" = "System, dnasr281@gmail.com"
    "G57" = "System, dnasr281@gmail.com"
    "G58" = "backup@backdoor.com, System"
    "G60" = "System, dnasr281@gmail.com"
    "G61" = "System, admin@admin.com"
    "G64" = "System, dnasr281@gmail.com"
    "G66" = "System, dnasr281@gmail.com"
    "G67" = "System, dnasr281@gmail.com"
    "G68" = "System, dnasr281@gmail.com"
    "G69" = "System, dnasr281@gmail.com"
    "G72" = "System, dnasr281@gmail.com"
    "G73" = "System, dnasr281@gmail.com"
    "G74" = "System, dnasr281@gmail.com"
    "G75" = "System, dnasr281@gmail.com"
    "G76" = "System, dnasr281@gmail.com"
    "G78" = "System, dnasr281@gmail.com"
    "G86" = "System, dnasr281@gmail.com"
    "G87" = "System, dnasr281@gmail.com"
    "G88" = "System, dnasr281@gmail.com"
    "G89" = "System, dnasr281@gmail.com"
    "G90" = "dnasr281@gmail.com, admin@admin.com"
    "G93" = "System, dnasr281@gmail.com"
    "G95" = "System, dnasr281@gmail.com"
    "G102" = "System, dnasr281@gmail.com"
    "G104" = "System, dnasr281@gmail.com"
    "G112" = "System, dnasr281@gmail.com"
    "G113" = "System, dnasr281@gmail.com"
    "G114" = "System, dnasr281@gmail.com"
    "G115" = "System, dnasr281@gmail.com"
    "G116" = "dnasr281@gmail.com, admin@admin.com"
    "G119" = "System, dnasr281@gmail.com"
    "G121" = "System, dnasr281@gmail.com"
    "G128" = "System, dnasr281@gmail.com"
    "G130" = "System, dnasr281@gmail.com"
    "G138" = "System, dnasr281@gmail.com"
    "G139" = "System, dnasr281@gmail.com"
    "G140" = "System, dnasr281@gmail.com"
    "G141" = "System, dnasr281@gmail.com"
    "G142" = "dnasr281@gmail.com, admin@admin.com"
    "G145" = "System, dnasr281@gmail.com"
    "G147" = "System, dnasr281@gmail.com"
    "G154" = "System, dnasr281@gmail.com"
    "G156" = "System, dnasr281@gmail.com"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
